# Niveaux.xlsx edit script
# 1) Fix the "manivelle" bug on sheet "Tous_les_niveaux": AA24 should contain the
#    same "M" marker as the rest of the Niveau 9 block (was left blank).
# 2) Add a second "Niveau 9" configuration block to sheet "Niveau 6-10" (rows 41-54),
#    the first copy being identical to the (now-fixed) block already present on
#    "Tous_les_niveaux", the second being a variant of it ("Configuration 2").
# 3) Update sheet/selection bookkeeping (active tab, active sheet, selections) to
#    match where the author ended up after making the change.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Niveau 1-5")
$ws2 = $wb.Worksheets.Item("Niveau 6-10")
$ws3 = $wb.Worksheets.Item("Tous_les_niveaux")

# ---------------------------------------------------------------------------
# 1) Bug fix on "Tous_les_niveaux": AA24 was blank, should show the "M" marker
#    (matches the pattern of the rest of this level-9 configuration block).
# ---------------------------------------------------------------------------
$ws3.Range("AA24").Value = "M"

# ---------------------------------------------------------------------------
# 2) New "Niveau 9" block (rows 41-46) on "Niveau 6-10", copied from the
#    already-existing (and now fixed) block on "Tous_les_niveaux" (Z21:AF26).
# ---------------------------------------------------------------------------

# Title row (merged C41:I41), formatting copied from the Niveau 6 title (row 3)
$ws2.Range("C3:H3").Copy($ws2.Range("C41:H41"))
$ws2.Range("K3").Copy($ws2.Range("I41"))
$ws2.Range("C41:I41").Merge()
$ws2.Range("C41").Value = "Niveau 9"

# Data rows, copied straight from the matching block on "Tous_les_niveaux"
$ws3.Range("Z22:AF22").Copy($ws2.Range("C42:I42"))
$ws3.Range("Z23:AF23").Copy($ws2.Range("C43:I43"))
$ws3.Range("Z24:AF24").Copy($ws2.Range("C44:I44"))
$ws3.Range("Z25:AF25").Copy($ws2.Range("C45:I45"))
$ws3.Range("Z26:AF26").Copy($ws2.Range("C46:I46"))

# ---------------------------------------------------------------------------
# 3) "Configuration 2 du niveau 9" block (rows 49-54) on "Niveau 6-10" - a
#    variant of the block above.
# ---------------------------------------------------------------------------

# Title row (merged C49:I49)
$ws2.Range("C3:H3").Copy($ws2.Range("C49:H49"))
$ws2.Range("K3").Copy($ws2.Range("I49"))
$ws2.Range("C49:I49").Merge()
$ws2.Range("C49").Value = "Configuration 2 du niveau 9"

# Start from the same data as the first block, then adjust the handful of
# cells that differ in this variant.
$ws2.Range("C42:I42").Copy($ws2.Range("C50:I50"))
$ws2.Range("C43:I43").Copy($ws2.Range("C51:I51"))
$ws2.Range("C44:I44").Copy($ws2.Range("C52:I52"))
$ws2.Range("C45:I45").Copy($ws2.Range("C53:I53"))
$ws2.Range("C46:I46").Copy($ws2.Range("C54:I54"))

Write-Host "done building blocks"
